# Applies the "CryCompanywiseStockReport_1" workbook edit described by the
# commit diff: a batch of stock-quantity / value corrections scattered
# throughout the sheet, plus a structural change near the end of the report
# (one line item row removed from the "XO FOOTWEAR PVT LTD" section, causing
# everything below it to shift up by one row, with subtotals/grand total
# recomputed accordingly).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Step 1: Delete row 926 ("XO-M Sport Shoe Rainbow (2 Colour)") entirely.
# This shifts every row below it up by one (926 downward), exactly matching
# the row-renumbering seen in the diff for rows 927-943 (before) -> 926-942
# (after), and updates the sheet's used range from M943 to M942.
# ---------------------------------------------------------------------------
$ws.Rows("926:926").Delete()

# ---------------------------------------------------------------------------
# Step 2: Apply every individual cell value correction - both the plain
# quantity/value corrections scattered across the sheet, and the follow-on
# fixups needed after the row-926 deletion (renumbering the "Sl No" column
# for the rows that shifted up, and recomputing the two subtotals + the
# grand total that are affected by the removed line item).
# ---------------------------------------------------------------------------
$changes = @(
    @{Cell="F22"; Value=11},
    @{Cell="G22"; Value=338.14},
    @{Cell="F35"; Value=186},
    @{Cell="G35"; Value=5241.48},
    @{Cell="B40"; Value=71755.91},
    @{Cell="F60"; Value=48},
    @{Cell="G60"; Value=787.6799999999999},
    @{Cell="F62"; Value=149},
    @{Cell="G62"; Value=8308.24},
    @{Cell="B73"; Value=258724.75},
    @{Cell="F230"; Value=10},
    @{Cell="G230"; Value=507.4},
    @{Cell="F233"; Value=76},
    @{Cell="G233"; Value=6232},
    @{Cell="F236"; Value=60},
    @{Cell="G236"; Value=2596.8},
    @{Cell="B247"; Value=88957.99000000001},
    @{Cell="F285"; Value=65},
    @{Cell="G285"; Value=1392.95},
    @{Cell="B291"; Value=52763.38},
    @{Cell="F401"; Value=22},
    @{Cell="G401"; Value=2286.02},
    @{Cell="F405"; Value=153},
    @{Cell="G405"; Value=26213.49},
    @{Cell="B409"; Value=31990.13},
    @{Cell="F475"; Value=55},
    @{Cell="G475"; Value=53339.55},
    @{Cell="B476"; Value=53339.55},
    @{Cell="F496"; Value=145},
    @{Cell="G496"; Value=3694.6},
    @{Cell="F506"; Value=89},
    @{Cell="G506"; Value=9155.43},
    @{Cell="F510"; Value=88},
    @{Cell="G510"; Value=21215.04},
    @{Cell="F517"; Value=71},
    @{Cell="G517"; Value=3883.7},
    @{Cell="B519"; Value=196578.66},
    @{Cell="B548"; Value=65068},
    @{Cell="E548"; Value=13.97},
    @{Cell="F548"; Value=0},
    @{Cell="G548"; Value=0},
    @{Cell="B549"; Value=53602},
    @{Cell="E549"; Value=15.69},
    @{Cell="F549"; Value=-232},
    @{Cell="G549"; Value=-3050.8},
    @{Cell="B550"; Value=65066},
    @{Cell="E550"; Value=13.61},
    @{Cell="F550"; Value=0},
    @{Cell="G550"; Value=0},
    @{Cell="B551"; Value=53263},
    @{Cell="E551"; Value=15.29},
    @{Cell="F551"; Value=-313},
    @{Cell="G551"; Value=-4009.53},
    @{Cell="B556"; Value=64922},
    @{Cell="E556"; Value=20.98},
    @{Cell="F556"; Value=0},
    @{Cell="G556"; Value=0},
    @{Cell="B557"; Value=45706},
    @{Cell="E557"; Value=23.58},
    @{Cell="F557"; Value=-207},
    @{Cell="G557"; Value=-4084.11},
    @{Cell="F580"; Value=169},
    @{Cell="G580"; Value=5882.89},
    @{Cell="F581"; Value=38},
    @{Cell="G581"; Value=2646.7},
    @{Cell="B584"; Value=35987.79},
    @{Cell="B640"; Value=53319},
    @{Cell="E640"; Value=310.64},
    @{Cell="F640"; Value=-6},
    @{Cell="G640"; Value=-1643.52},
    @{Cell="B641"; Value=64810},
    @{Cell="E641"; Value=291.22},
    @{Cell="F641"; Value=2},
    @{Cell="G641"; Value=547.84},
    @{Cell="B669"; Value=60022},
    @{Cell="E669"; Value=37.22},
    @{Cell="F669"; Value=-113},
    @{Cell="G669"; Value=-3709.79},
    @{Cell="B670"; Value=64830},
    @{Cell="E670"; Value=34.9},
    @{Cell="F670"; Value=89},
    @{Cell="G670"; Value=2921.87},
    @{Cell="F691"; Value=35},
    @{Cell="G691"; Value=3091.9},
    @{Cell="B692"; Value=160287.05},
    @{Cell="F718"; Value=0},
    @{Cell="G718"; Value=0},
    @{Cell="B733"; Value=83201.81},
    @{Cell="F800"; Value=230},
    @{Cell="G800"; Value=30613},
    @{Cell="B803"; Value=31423.06},
    @{Cell="F818"; Value=51},
    @{Cell="G818"; Value=7281.27},
    @{Cell="B839"; Value=278053.05},
    @{Cell="F875"; Value=75},
    @{Cell="G875"; Value=4015.5},
    @{Cell="F882"; Value=3},
    @{Cell="G882"; Value=123.15},
    @{Cell="B884"; Value=20113.66},
    @{Cell="F889"; Value=92},
    @{Cell="G889"; Value=2781.16},
    @{Cell="F890"; Value=1488},
    @{Cell="G890"; Value=242707.68},
    @{Cell="F893"; Value=54},
    @{Cell="G893"; Value=7811.1},
    @{Cell="B896"; Value=268766.54},
    @{Cell="F908"; Value=22},
    @{Cell="G908"; Value=3501.52},
    @{Cell="B912"; Value=16302.98},
    @{Cell="F922"; Value=20},
    @{Cell="G922"; Value=11162.6},
    @{Cell="A926"; Value=766},
    @{Cell="A927"; Value=767},
    @{Cell="A928"; Value=768},
    @{Cell="A929"; Value=769},
    @{Cell="A930"; Value=770},
    @{Cell="A931"; Value=771},
    @{Cell="A932"; Value=772},
    @{Cell="A933"; Value=773},
    @{Cell="B934"; Value=90830.57},
    @{Cell="A936"; Value=774},
    @{Cell="A937"; Value=775},
    @{Cell="A938"; Value=776},
    @{Cell="B940"; Value=3941036.75},
    @{Cell="B941"; Value=3941036.75}
)

foreach ($change in $changes) {
    $ws.Range($change.Cell).Value2 = $change.Value
}
